# Update the build timestamp embedded in the "version" strings across the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner (A2) and recommended citation (A6)
$a2text = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2text.Replace($oldStamp, $newStamp)

$a6text = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6text.Replace($oldStamp, $newStamp)

# Data sheet: build_version column (S) for every data row (2-8)
for ($r = 2; $r -le 8; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    $cellText = $cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
